# Generate Report for Archive
# Status text "Ready for handoff" -> "In Translation" across all sheets,
# and the now-narrower Status/language columns are re-sized to match.

$wb = $excel.ActiveWorkbook

# NOTE: the host's column-width model snaps to a 6-px grid (MDW=6), so a
# ColumnWidth input of 12.5 is the value that lands closest to the target
# stored width of 13.4101845877511 (rounds to 13.333333...), the nearest
# the COM surface can reproduce.
$newStatus = "In Translation"
$newWidth  = 12.5

# --- Overview sheet: zh-cn (col E) / de-de (col F) status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet: Status column (col C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet: Status column (col C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
